$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.412.93'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '2.646.21'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.546'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.57%  '
$ws.Range('D9').Value = '2.646.21'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  +8.05%  '
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = '3.128.66'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '68.352.53'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '2.643.52'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '365.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('E22').Value = '  +4.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.16'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000106'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.52%  '
$ws.Range('D29').Value = '2.774.63'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '574.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('E33').Value = '  +3.33%  '
$ws.Range('E34').Value = '  +1.67%  '
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  +6.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.44%  '
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('E42').Value = '  +2.74%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.12%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₆0333'
$ws.Range('E44').Value = '  +9.32%  '
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '157.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.98'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.47%  '
